$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.984.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.52%  '
$ws.Range('E2').Style = 'Normal'

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.656.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +10.72%  '
$ws.Range('E3').Style = 'Normal'

# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('E4').Style = 'Normal'

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.76%  '
$ws.Range('E5').Style = 'Normal'

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +7.18%  '
$ws.Range('E6').Style = 'Normal'

# Row 7: XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.612'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +8.57%  '
$ws.Range('E7').Style = 'Normal'

# Row 8: USDC
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('E8').Style = 'Normal'

# Row 9: Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.599'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +16.99%  '
$ws.Range('E9').Style = 'Normal'

# Row 10: Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +13.72%  '
$ws.Range('E10').Style = 'Normal'

# Row 11: OKB
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.26'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.57%  '
$ws.Range('E11').Style = 'Normal'

# Row 12: Dogecoin
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0857'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.30%  '
$ws.Range('E12').Style = 'Normal'

# Row 13: Polkadot
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +17.66%  '
$ws.Range('E13').Style = 'Normal'

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.057.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +10.27%  '
$ws.Range('E14').Style = 'Normal'

# Row 15: TRON
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('E15').Style = 'Normal'

# Row 16: WrappedEther
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.660.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +10.46%  '
$ws.Range('E16').Style = 'Normal'

# Row 17: Polygon
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.948'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +12.35%  '
$ws.Range('E17').Style = 'Normal'

# Row 18: Chainlink
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +7.60%  '
$ws.Range('E18').Style = 'Normal'

# Row 19: WrappedBTC
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '47.542.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('E19').Style = 'Normal'

# Row 20: ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000104'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.93%  '
$ws.Range('E20').Style = 'Normal'

# Row 21: InternetComputer(DFINITY)
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.23%  '
$ws.Range('E21').Style = 'Normal'

# Row 22: Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +10.39%  '
$ws.Range('E22').Style = 'Normal'

# Row 23: Litecoin
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.73%  '
$ws.Range('E23').Style = 'Normal'

# Row 24: BitcoinCash
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +10.33%  '
$ws.Range('E24').Style = 'Normal'

# Row 25: PancakeSwap
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +10.42%  '
$ws.Range('E25').Style = 'Normal'

# Row 26: EthereumClassic
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '31.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +49.40%  '
$ws.Range('E26').Style = 'Normal'

# Row 27: ImmutableX
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +16.69%  '
$ws.Range('E27').Style = 'Normal'

# Row 28: Dai
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E28').Style = 'Normal'

# Row 29: LEO
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('E29').Style = 'Normal'

# Row 30: Cosmos
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +10.57%  '
$ws.Range('E30').Style = 'Normal'

# Row 31: InjectiveProtocol
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '41.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.76%  '
$ws.Range('E31').Style = 'Normal'

# Row 32: Toncoin
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.24%  '
$ws.Range('E32').Style = 'Normal'

# Row 33: Filecoin
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +13.86%  '
$ws.Range('E33').Style = 'Normal'

# Row 34: LidoDAOToken
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('E34').Style = 'Normal'

# Row 35: ARBITRUM
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +18.42%  '
$ws.Range('E35').Style = 'Normal'

# Row 36: WEMIXToken
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.86%  '
$ws.Range('E36').Style = 'Normal'

# Row 37: Hedera
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +9.63%  '
$ws.Range('E37').Style = 'Normal'

# Row 38: Monero
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '153.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('E38').Style = 'Normal'

# Row 39: Kaspa
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.79%  '
$ws.Range('E39').Style = 'Normal'

# Row 40: Stellar
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.124'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.76%  '
$ws.Range('E40').Style = 'Normal'

# Row 41: Celestia
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.03'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +12.21%  '
$ws.Range('E41').Style = 'Normal'

# Row 42: RenderToken
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +12.73%  '
$ws.Range('E42').Style = 'Normal'

# Row 43: NEARProtocol
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +15.92%  '
$ws.Range('E43').Style = 'Normal'

# Row 44: EnergySwap
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.72'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +44.86%  '
$ws.Range('E44').Style = 'Normal'

# Row 45: VeChain
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0333'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +10.56%  '
$ws.Range('E45').Style = 'Normal'

# Row 46: Maker
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.095.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +7.87%  '
$ws.Range('E46').Style = 'Normal'

# Row 47: FirstDigitalUSD
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E47').Style = 'Normal'

# Row 48: BitcoinSV
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.06%  '
$ws.Range('E48').Style = 'Normal'

# Row 49: Aave
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '115.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +14.32%  '
$ws.Range('E49').Style = 'Normal'

# Row 50: Stacks
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.35%  '
$ws.Range('E50').Style = 'Normal'

# Row 51: FraxShare
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.90%  '
$ws.Range('E51').Style = 'Normal'
